$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "34.623.14"
$ws.Cells.Item(2, 5).Value = "  +1.72%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.804.19"
$ws.Cells.Item(3, 5).Value = "  +0.95%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.16%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "224.56"
$ws.Cells.Item(5, 5).Value = "  -1.37%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "0.553"
$ws.Cells.Item(6, 5).Value = "  -0.03%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.12%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "32.56"
$ws.Cells.Item(8, 5).Value = "  +3.43%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.290"
$ws.Cells.Item(9, 5).Value = "  +3.17%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "0.0710"
$ws.Cells.Item(10, 5).Value = "  +7.59%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "0.0929"
$ws.Cells.Item(11, 5).Value = "  +0.08%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "2.063.79"
$ws.Cells.Item(12, 5).Value = "  +0.89%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -3.58%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "1.804.87"
$ws.Cells.Item(14, 5).Value = "  +0.91%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  +1.17%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "34.633.65"
$ws.Cells.Item(16, 5).Value = "  +1.60%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "4.33"

# Row 18
$ws.Cells.Item(18, 4).Value = "69.23"
$ws.Cells.Item(18, 5).Value = "  -0.34%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "253.28"
$ws.Cells.Item(19, 5).Value = "  +0.08%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +7.99%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "11.10"
$ws.Cells.Item(21, 5).Value = "  +5.84%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -0.17%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "4.25"
$ws.Cells.Item(23, 5).Value = "  -0.44%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +1.39%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "161.56"
$ws.Cells.Item(25, 5).Value = "  +2.69%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "16.43"
$ws.Cells.Item(26, 5).Value = "  -0.89%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +1.72%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -0.13%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -0.17%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "Swop.fi"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/yrCr2HW2c+swopfi-swop"
$ws.Cells.Item(30, 4).Value = "578.95"
$ws.Cells.Item(30, 5).Value = "  +1,002.65%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "Hedera"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(31, 4).Value = "0.0529"
$ws.Cells.Item(31, 5).Value = "  +2.39%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Filecoin"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(32, 4).Value = "3.80"
$ws.Cells.Item(32, 5).Value = "  -0.25%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "PancakeSwap"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(33, 4).Value = "1.20"
$ws.Cells.Item(33, 5).Value = "  -0.56%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(34, 4).Value = "3.63"
$ws.Cells.Item(34, 5).Value = "  +0.61%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "LidoDAOToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(35, 4).Value = "1.89"
$ws.Cells.Item(35, 5).Value = "  +2.11%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "Maker"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(36, 4).Value = "1.434.47"
$ws.Cells.Item(36, 5).Value = "  -1.28%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "TrustWalletToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(37, 4).Value = "1.07"
$ws.Cells.Item(37, 5).Value = "  +0.11%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "ImmutableX"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(38, 4).Value = "0.644"
$ws.Cells.Item(38, 5).Value = "  +2.27%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "VeChain"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(39, 4).Value = "0.0192"
$ws.Cells.Item(39, 5).Value = "  +2.85%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Aave"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(40, 4).Value = "84.90"
$ws.Cells.Item(40, 5).Value = "  +1.81%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "ARBITRUM"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(41, 4).Value = "0.957"
$ws.Cells.Item(41, 5).Value = "  +6.08%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "MXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(42, 4).Value = "2.80"
$ws.Cells.Item(42, 5).Value = "  -0.83%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "HuobiToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(43, 4).Value = "2.35"
$ws.Cells.Item(43, 5).Value = "  -0.03%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "RenderToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(44, 4).Value = "2.16"
$ws.Cells.Item(44, 5).Value = "  +3.85%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "FraxShare"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(45, 4).Value = "6.04"
$ws.Cells.Item(45, 5).Value = "  +4.40%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "WEMIXToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(46, 4).Value = "1.06"
$ws.Cells.Item(46, 5).Value = "  -1.06%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Kaspa"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(47, 4).Value = "0.0498"
$ws.Cells.Item(47, 5).Value = "  -2.69%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "RocketPoolETH"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(48, 4).Value = "1.958.57"
$ws.Cells.Item(48, 5).Value = "  +0.52%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "106.47"
$ws.Cells.Item(49, 5).Value = "  +8.67%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(50, 4).Value = "12.27"
$ws.Cells.Item(50, 5).Value = "  +2.50%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "PaxDollar"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(51, 4).Value = "1.00"
$ws.Cells.Item(51, 5).Value = "  -0.04%  "
